# Add a new cage record (search/info form entry) to the bottom of the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 32
$ws.Cells.Item($newRow, 1).Value = "22A"
$ws.Cells.Item($newRow, 2).Value = 23
$ws.Cells.Item($newRow, 3).Value = 23
$ws.Cells.Item($newRow, 4).Value = 22
$ws.Cells.Item($newRow, 5).Value = "Plastic"
